$p = $ppt.ActivePresentation

# Remove the slide with SlideID 262 (the last slide, slide7.xml) as per the
# commit's change: the screen-definition slide deck drops its final slide.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 262) {
        $s.Delete()
    }
}
